$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.070.47'
$ws.Range('D2').Style = $origStyle
$origStyle = $ws.Range('E2').Style
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('E2').Style = $origStyle
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.178.25'
$ws.Range('D3').Style = $origStyle
$origStyle = $ws.Range('E3').Style
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -4.02%  '
$ws.Range('E3').Style = $origStyle
$origStyle = $ws.Range('E4').Style
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E4').Style = $origStyle
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.94'
$ws.Range('D5').Style = $origStyle
$origStyle = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.43%  '
$ws.Range('E5').Style = $origStyle
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.04'
$ws.Range('D6').Style = $origStyle
$origStyle = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.77%  '
$ws.Range('E6').Style = $origStyle
$origStyle = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E7').Style = $origStyle
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.174.98'
$ws.Range('D8').Style = $origStyle
$origStyle = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('E8').Style = $origStyle
$origStyle = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('E9').Style = $origStyle
$origStyle = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -6.28%  '
$ws.Range('E10').Style = $origStyle
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.21'
$ws.Range('D11').Style = $origStyle
$origStyle = $ws.Range('E11').Style
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.87%  '
$ws.Range('E11').Style = $origStyle
$origStyle = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.53%  '
$ws.Range('E12').Style = $origStyle
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000236'
$ws.Range('D13').Style = $origStyle
$origStyle = $ws.Range('E13').Style
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.67%  '
$ws.Range('E13').Style = $origStyle
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.52'
$ws.Range('D14').Style = $origStyle
$origStyle = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('E14').Style = $origStyle
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.692.35'
$ws.Range('D15').Style = $origStyle
$origStyle = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.28%  '
$ws.Range('E15').Style = $origStyle
$origStyle = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('E16').Style = $origStyle
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.180.39'
$ws.Range('D17').Style = $origStyle
$origStyle = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.97%  '
$ws.Range('E17').Style = $origStyle
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.020.08'
$ws.Range('D18').Style = $origStyle
$origStyle = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('E18').Style = $origStyle
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.55'
$ws.Range('D19').Style = $origStyle
$origStyle = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -4.57%  '
$ws.Range('E19').Style = $origStyle
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.41'
$ws.Range('D20').Style = $origStyle
$origStyle = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.48%  '
$ws.Range('E20').Style = $origStyle
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.01'
$ws.Range('D21').Style = $origStyle
$origStyle = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('E21').Style = $origStyle
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.697'
$ws.Range('D22').Style = $origStyle
$origStyle = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.00%  '
$ws.Range('E22').Style = $origStyle
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.60'
$ws.Range('D23').Style = $origStyle
$origStyle = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.79%  '
$ws.Range('E23').Style = $origStyle
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.32'
$ws.Range('D24').Style = $origStyle
$origStyle = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.72%  '
$ws.Range('E24').Style = $origStyle
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.40'
$ws.Range('D25').Style = $origStyle
$origStyle = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.35%  '
$ws.Range('E25').Style = $origStyle
$origStyle = $ws.Range('E27').Style
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E27').Style = $origStyle
$origStyle = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.00%  '
$ws.Range('E28').Style = $origStyle
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.74'
$ws.Range('D29').Style = $origStyle
$origStyle = $ws.Range('E29').Style
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.02%  '
$ws.Range('E29').Style = $origStyle
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.65'
$ws.Range('D30').Style = $origStyle
$origStyle = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.04%  '
$ws.Range('E30').Style = $origStyle
$origStyle = $ws.Range('E31').Style
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.62%  '
$ws.Range('E31').Style = $origStyle
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.18'
$ws.Range('D32').Style = $origStyle
$origStyle = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.92%  '
$ws.Range('E32').Style = $origStyle
$origStyle = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.96%  '
$ws.Range('E33').Style = $origStyle
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.36'
$ws.Range('D34').Style = $origStyle
$origStyle = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.26%  '
$ws.Range('E34').Style = $origStyle
$origStyle = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.95%  '
$ws.Range('E35').Style = $origStyle
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.80'
$ws.Range('D36').Style = $origStyle
$origStyle = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.47%  '
$ws.Range('E36').Style = $origStyle
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.34'
$ws.Range('D37').Style = $origStyle
$origStyle = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.18%  '
$ws.Range('E37').Style = $origStyle
$origStyle = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.15%  '
$ws.Range('E38').Style = $origStyle
$origStyle = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.20%  '
$ws.Range('E39').Style = $origStyle
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '404.99'
$ws.Range('D40').Style = $origStyle
$origStyle = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -6.90%  '
$ws.Range('E40').Style = $origStyle
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.07'
$ws.Range('D41').Style = $origStyle
$origStyle = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('E41').Style = $origStyle
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.64'
$ws.Range('D42').Style = $origStyle
$origStyle = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('E42').Style = $origStyle
$origStyle = $ws.Range('E43').Style
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.94%  '
$ws.Range('E43').Style = $origStyle
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.814.40'
$ws.Range('D44').Style = $origStyle
$origStyle = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -9.60%  '
$ws.Range('E44').Style = $origStyle
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.251'
$ws.Range('D45').Style = $origStyle
$origStyle = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.81%  '
$ws.Range('E45').Style = $origStyle
$origStyle = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.71%  '
$ws.Range('E47').Style = $origStyle
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.31'
$ws.Range('D48').Style = $origStyle
$origStyle = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('E48').Style = $origStyle
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.00'
$ws.Range('D49').Style = $origStyle
$origStyle = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.02%  '
$ws.Range('E49').Style = $origStyle
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.08'
$ws.Range('D50').Style = $origStyle
$origStyle = $ws.Range('E50').Style
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.83%  '
$ws.Range('E50').Style = $origStyle
$origStyle = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.07%  '
$ws.Range('E51').Style = $origStyle
